$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column J (2022) into the new column K (2023)
$ws.Range("J3:J6").Copy()
$ws.Range("K3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new 2023 column values
$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 1736.9
$ws.Range("K5").Value = 1239.3
$ws.Range("K6").Value = 1934.5
